# Actualización automática 2025-09-11 09:25:09
#
# A new advisor record ("JUNCO SANCHEZ ARTURO ENRIQUE", under
# "OFICINA-CATAECSA") was added to the source data, alphabetically placed
# right before "KITCHENSCO S.A.". This inserts one row at that position in
# both the "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets, shifting every
# following row down by one (all its values are 0, matching the blank
# template rows around it). The trailing "N de <count>" summary row on
# "VENTAS POR GRUPO" is updated to reflect the new total row count
# (317 -> 318).

$wb = $excel.ActiveWorkbook

# ----- Sheet "VENTAS POR GRUPO": insert new row at 279 -----
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Rows.Item(279).Insert()
$ws1.Range("A279").Value = "OFICINA-CATAECSA"
$ws1.Range("B279").Value = "JUNCO SANCHEZ ARTURO ENRIQUE"
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(279, $c).Value = 0
}

# Update the "X de 317" -> "X de 318" footer/summary row, which used to be
# row 319 and is now row 320 after the insert above.
$summaryRow1 = 320
for ($c = 3; $c -le 18; $c++) {
    $cell = $ws1.Cells.Item($summaryRow1, $c)
    $cell.Value = $cell.Text.Replace("de 317", "de 318")
}

# ----- Sheet "VENTA MENSUAL": insert new row at 283 -----
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Rows.Item(283).Insert()
$ws2.Range("A283").Value = "OFICINA-CATAECSA"
$ws2.Range("B283").Value = "JUNCO SANCHEZ ARTURO ENRIQUE"
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(283, $c).Value = 0
}
